# Insert a new row before row 93, shifting existing rows 93..191 down to 94..192,
# then populate the new row 93 with the latest price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 93 (pushes old row 93 -> 94, ..., old row 191 -> 192)
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new data record
$ws.Cells.Item(93, 1).Value = 4
$ws.Cells.Item(93, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value = "Los Lagos"
$ws.Cells.Item(93, 4).Value = 44539
$ws.Cells.Item(93, 5).Value = 10
$ws.Cells.Item(93, 6).Value = 100112043
$ws.Cells.Item(93, 7).Value = "Pepino ensalada"
$ws.Cells.Item(93, 8).Value = "Sin especificar"
$ws.Cells.Item(93, 9).Value = "Primera"
$ws.Cells.Item(93, 10).Value = 150
$ws.Cells.Item(93, 11).Value = 10000
$ws.Cells.Item(93, 12).Value = 11000
$ws.Cells.Item(93, 13).Value = 10533
$ws.Cells.Item(93, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(93, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(93, 16).Value = 176
$ws.Cells.Item(93, 17).Value = 60
$ws.Cells.Item(93, 18).Value = "Hortaliza"
